$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 188, shifting existing rows 188:198 down to 189:199
$ws.Rows.Item(188).Insert()

# Populate the new row 188 with its values (same template as the other Apio rows,
# only Fecha and Volumen differ; price fields match the old row 188's values)
$ws.Cells.Item(188, 1).Value = 5
$ws.Cells.Item(188, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(188, 3).Value = "Maule"
$ws.Cells.Item(188, 4).Value = 44753
$ws.Cells.Item(188, 5).Value = 7
$ws.Cells.Item(188, 6).Value = 100112017
$ws.Cells.Item(188, 7).Value = "Apio"
$ws.Cells.Item(188, 8).Value = "Americana (o)"
$ws.Cells.Item(188, 9).Value = "Primera"
$ws.Cells.Item(188, 10).Value = 700
$ws.Cells.Item(188, 11).Value = 7000
$ws.Cells.Item(188, 12).Value = 7000
$ws.Cells.Item(188, 13).Value = 7000
$ws.Cells.Item(188, 14).Value = "$/docena de matas"
$ws.Cells.Item(188, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(188, 16).Value = 1167
$ws.Cells.Item(188, 17).Value = 6
$ws.Cells.Item(188, 18).Value = "Hortaliza"
